# "matrice and a few details" - update implication-matrix percentages and
# give the still-empty row 10 (task "10") the same percentage formatting as
# the rows above it, then leave the view scrolled/selected where the author
# left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the implication-matrix percentages (rows 10-14, columns D/E) ---
$ws.Range("D10").Value = 0.3
$ws.Range("E10").Value = 0.7

$ws.Range("D11").Value = 0.7
$ws.Range("E11").Value = 0.3

$ws.Range("D12").Value = 0.7
$ws.Range("E12").Value = 0.3

$ws.Range("D13").Value = 0.7
$ws.Range("E13").Value = 0.3

$ws.Range("D14").Value = 0.9
$ws.Range("E14").Value = 0.1

# --- Give row 16 (task n°10) the same percentage number format used by the
#     rows above (D7:E15), even though the cells themselves stay empty ---
$ws.Range("D16:E16").NumberFormat = "0%"

# --- Restore the view state (scroll position + active selection) left by
#     the author when the workbook was saved ---
$ws.Range("J14").Select()
$excel.ActiveWindow.ScrollRow = 3
